$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) serial date values from 45233 to 45243 for rows 2-32
for ($row = 2; $row -le 32; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 45243
}
